# Update "想去人数" (F column) counts across sheets to reflect freshly
# scraped data (re-generated gh-pages output).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 285
$ws.Range("F4").Value = 641
$ws.Range("F5").Value = 2855
$ws.Range("F9").Value = 284
$ws.Range("F11").Value = 5120
$ws.Range("F12").Value = 16
$ws.Range("F22").Value = 2145
$ws.Range("F23").Value = 1375
$ws.Range("F25").Value = 70
$ws.Range("F26").Value = 1028
$ws.Range("F30").Value = 1561
$ws.Range("F32").Value = 15
$ws.Range("F34").Value = 1127
$ws.Range("F36").Value = 559
$ws.Range("F38").Value = 355
$ws.Range("F41").Value = 64
$ws.Range("F44").Value = 8
$ws.Range("F46").Value = 107
$ws.Range("F48").Value = 116
$ws.Range("F49").Value = 399

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 115

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 1590
$ws.Range("F9").Value = 2631
$ws.Range("F10").Value = 927

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 641
$ws.Range("F6").Value = 2855
$ws.Range("F7").Value = 1590
$ws.Range("F8").Value = 284
$ws.Range("F9").Value = 2631
$ws.Range("F11").Value = 927
$ws.Range("F13").Value = 5120
$ws.Range("F19").Value = 115
$ws.Range("F23").Value = 2145
$ws.Range("F24").Value = 1375
$ws.Range("F26").Value = 70
$ws.Range("F28").Value = 1028
$ws.Range("F31").Value = 1561
$ws.Range("F33").Value = 15
$ws.Range("F35").Value = 1127
$ws.Range("F36").Value = 559
$ws.Range("F39").Value = 355
$ws.Range("F46").Value = 107
$ws.Range("F48").Value = 399

$wb.Save()
